# Scheduled runner update: refresh cached market-price figures (and the
# derived profit columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1666.6
$ws.Range("I32").Value = 716.5
$ws.Range("K32").Value = 716.5
$ws.Range("M32").Value = -390.5

$ws.Range("H40").Value = 1903.3334
$ws.Range("J40").Value = 2400
$ws.Range("L40").Value = 2400
$ws.Range("N40").Value = -2750

$ws.Range("H96").Value = 897.3333
$ws.Range("I96").Value = 860.8570999999999
$ws.Range("J96").Value = 1025
$ws.Range("K96").Value = 2582.5713
$ws.Range("L96").Value = 3075
$ws.Range("M96").Value = -1209.5713
$ws.Range("N96").Value = -5821

$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = 18000
$ws.Range("N125").Value = -22920

$ws.Range("H128").Value = 48000
$ws.Range("J128").Value = 48000
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960

$ws.Range("H129").Value = 1429101.4
$ws.Range("J129").Value = 2500652.5
$ws.Range("L129").Value = 7501957.5
$ws.Range("N129").Value = -7511957.5

$ws.Range("H131").Value = 3398.3333
$ws.Range("I131").Value = 2697.5
$ws.Range("K131").Value = 8092.5
$ws.Range("M131").Value = -3052.5

$ws.Range("H138").Value = 1936.0513
$ws.Range("J138").Value = 2330.4375
$ws.Range("L138").Value = 6991.3125
$ws.Range("N138").Value = -17271.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 66

$ws.Range("H32").Value = 26189.209
$ws.Range("I32").Value = 27210.55
$ws.Range("J32").Value = 12571.333
$ws.Range("K32").Value = 27210.55
$ws.Range("L32").Value = 12571.333
$ws.Range("M32").Value = -26923.55
$ws.Range("N32").Value = -13145.333

$ws.Range("H61").Value = 2455.0344
$ws.Range("I61").Value = 1976.4546
$ws.Range("J61").Value = 3959.1428
$ws.Range("K61").Value = 1976.4546
$ws.Range("L61").Value = 3959.1428
$ws.Range("M61").Value = -1764.4546
$ws.Range("N61").Value = -4383.1428

$ws.Range("H74").Value = 50002710
$ws.Range("I74").Value = 52634340
$ws.Range("K74").Value = 52634340
$ws.Range("M74").Value = -52633466

$ws.Range("H77").Value = 50002710
$ws.Range("I77").Value = 52634340
$ws.Range("K77").Value = 263171700
$ws.Range("M77").Value = -263167332

$ws.Range("H102").Value = 1483.6
$ws.Range("I102").Value = 1216
$ws.Range("J102").Value = 1885
$ws.Range("K102").Value = 1216
$ws.Range("L102").Value = 1885
$ws.Range("M102").Value = 406
$ws.Range("N102").Value = -5129

$ws.Range("H132").Value = 12781.066
$ws.Range("I132").Value = 1463.1351
$ws.Range("J132").Value = 65126.5
$ws.Range("K132").Value = 4389.4053
$ws.Range("L132").Value = 195379.5
$ws.Range("M132").Value = -1859.4053
$ws.Range("N132").Value = -200439.5

$ws.Range("H136").Value = 2455.0344
$ws.Range("I136").Value = 1976.4546
$ws.Range("J136").Value = 3959.1428
$ws.Range("K136").Value = 5929.3638
$ws.Range("L136").Value = 11877.4284
$ws.Range("M136").Value = -3379.3638
$ws.Range("N136").Value = -16977.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3336766.2
$ws.Range("I7").Value = 5000149.5
$ws.Range("K7").Value = 5000149.5
$ws.Range("M7").Value = -5000036.5

$ws.Range("H39").Value = 13999
$ws.Range("I39").Value = 13999
$ws.Range("K39").Value = 13999
$ws.Range("M39").Value = -13610

$ws.Range("H88").Value = 23749.75
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 28333
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 28333
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -29145

$ws.Range("H91").Value = 23749.75
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 28333
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 28333
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -31141

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14574.56
$ws.Range("I31").Value = 19966.5
$ws.Range("K31").Value = 19966.5
$ws.Range("M31").Value = -19671.5

$ws.Range("H34").Value = 14574.56
$ws.Range("I34").Value = 19966.5
$ws.Range("K34").Value = 19966.5
$ws.Range("M34").Value = -19764.5

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0

$ws.Range("H95").Value = 14666.667
$ws.Range("J95").Value = 14666.667
$ws.Range("L95").Value = 14666.667
$ws.Range("N95").Value = -20158.667

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0

$ws.Range("H105").Value = 5953453.5
$ws.Range("I105").Value = 8929255
$ws.Range("J105").Value = 1851.4286
$ws.Range("K105").Value = 8929255
$ws.Range("L105").Value = 1851.4286
$ws.Range("M105").Value = -8927508
$ws.Range("N105").Value = -5345.4286

$ws.Range("H107").Value = 931.94116
$ws.Range("I107").Value = 233.44444
$ws.Range("J107").Value = 1717.75
$ws.Range("K107").Value = 233.44444
$ws.Range("L107").Value = 1717.75
$ws.Range("M107").Value = 1686.55556
$ws.Range("N107").Value = -5557.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 839.8570999999999
$ws.Range("I11").Value = 59.666668
$ws.Range("J11").Value = 1425
$ws.Range("K11").Value = 179.000004
$ws.Range("L11").Value = 4275
$ws.Range("M11").Value = -39.00000399999999
$ws.Range("N11").Value = -4555

$ws.Range("H36").Value = 114526.25
$ws.Range("I36").Value = 2401
$ws.Range("K36").Value = 7203
$ws.Range("M36").Value = -7034

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25002602
$ws.Range("I102").Value = 33336142
$ws.Range("J102").Value = 1986.8
$ws.Range("K102").Value = 33336142
$ws.Range("L102").Value = 1986.8
$ws.Range("M102").Value = -33334520
$ws.Range("N102").Value = -5230.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2563.7273
$ws.Range("J22").Value = 4433.3335
$ws.Range("L22").Value = 4433.3335
$ws.Range("N22").Value = -5023.3335

$ws.Range("H27").Value = 2563.7273
$ws.Range("J27").Value = 4433.3335
$ws.Range("L27").Value = 4433.3335
$ws.Range("N27").Value = -4647.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1552.8286
$ws.Range("I122").Value = 1618.5416
$ws.Range("J122").Value = 1409.4546
$ws.Range("K122").Value = 4855.6248
$ws.Range("L122").Value = 4228.3638
$ws.Range("M122").Value = -2405.6248
$ws.Range("N122").Value = -9128.363799999999

$ws.Range("H126").Value = 1316.88
$ws.Range("J126").Value = 1066.9166
$ws.Range("L126").Value = 3200.7498
$ws.Range("N126").Value = -8140.7498

$ws.Range("H136").Value = 27028930
$ws.Range("I136").Value = 41668396
$ws.Range("K136").Value = 125005188
$ws.Range("M136").Value = -125002638
